# Append new trading log rows (32-35) to the sheet, matching the source
# commit "Update trading results - Mon Sep 22 12:43:50 UTC 2025".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: TRADING_ATTEMPT XRP
$ws.Cells.Item(32, 1).Value = "2025-09-22T12:43:46.357167"
$ws.Cells.Item(32, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(32, 3).Value = "XRP"
$ws.Cells.Item(32, 4).Value = "UNKNOWN"
$ws.Cells.Item(32, 5).Value = 2.973349691321248
$ws.Cells.Item(32, 11).Value = "ATTEMPT"
$ws.Cells.Item(32, 12).Value = "Attempting trade 1/2"

# Row 33: POSITION_OPENED XRP
$ws.Cells.Item(33, 1).Value = "2025-09-22T12:43:48.480932"
$ws.Cells.Item(33, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(33, 3).Value = "XRP"
$ws.Cells.Item(33, 4).Value = "UNKNOWN"
$ws.Cells.Item(33, 5).Value = 2.973349691321248
$ws.Cells.Item(33, 6).Value = 2400
$ws.Cells.Item(33, 7).Value = 20
$ws.Cells.Item(33, 8).Value = 0.03925898863075103
$ws.Cells.Item(33, 11).Value = "SUCCESS"

# Row 34: TRADING_ATTEMPT ENA
$ws.Cells.Item(34, 1).Value = "2025-09-22T12:43:48.500959"
$ws.Cells.Item(34, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(34, 3).Value = "ENA"
$ws.Cells.Item(34, 4).Value = "UNKNOWN"
$ws.Cells.Item(34, 5).Value = 0.6448516400994989
$ws.Cells.Item(34, 11).Value = "ATTEMPT"
$ws.Cells.Item(34, 12).Value = "Attempting trade 2/2"

# Row 35: POSITION_OPENED ENA
$ws.Cells.Item(35, 1).Value = "2025-09-22T12:43:50.177958"
$ws.Cells.Item(35, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(35, 3).Value = "ENA"
$ws.Cells.Item(35, 4).Value = "UNKNOWN"
$ws.Cells.Item(35, 5).Value = 0.6448516400994989
$ws.Cells.Item(35, 6).Value = 2400
$ws.Cells.Item(35, 7).Value = 10
$ws.Cells.Item(35, 8).Value = 0.544849307081386
$ws.Cells.Item(35, 11).Value = "SUCCESS"
